$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

function Set-NumberCell($row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

$rows = @(
    @{
        Row = 71
        Date = "2025-09-06"
        Liga = "Liga de Expansión MX"
        Local = "Mineros de Zacatecas"
        Visitante = "Cancún"
        Prediccion = "Home Win"
        Probabilidad = "62.18%"
        Cuota = 1.95
        EV = "20.04%"
        Stake = 1.3
        StakePct = 0.02236911800293895
        KellyFrac = 0.2236911800293895
        Status = "Pending"
    },
    @{
        Row = 72
        Date = "2025-09-06"
        Liga = "Liga de Expansión MX"
        Local = "Dorados"
        Visitante = "Irapuato"
        Prediccion = "Away Win"
        Probabilidad = "49.64%"
        Cuota = 2.2
        EV = "8.12%"
        Stake = 0.5
        StakePct = 0.007673772179820727
        KellyFrac = 0.07673772179820727
        Status = "Pending"
    },
    @{
        Row = 73
        Date = "2025-09-06"
        Liga = "Liga de Expansión MX"
        Local = "Tlaxcala"
        Visitante = "Alebrijes de Oaxaca"
        Prediccion = "Home Win"
        Probabilidad = "72.66%"
        Cuota = 1.75
        EV = "25.88%"
        Stake = 2.1
        StakePct = 0.03620043620502014
        KellyFrac = 0.3620043620502013
        Status = "Pending"
    }
)

foreach ($r in $rows) {
    $row = $r.Row
    Set-TextCell $row 1 $r.Date
    Set-TextCell $row 2 $r.Liga
    Set-TextCell $row 3 $r.Local
    Set-TextCell $row 4 $r.Visitante
    Set-TextCell $row 5 $r.Prediccion
    Set-TextCell $row 6 $r.Probabilidad
    Set-NumberCell $row 7 $r.Cuota
    Set-TextCell $row 8 $r.EV
    Set-NumberCell $row 9 $r.Stake
    Set-NumberCell $row 10 $r.StakePct
    Set-NumberCell $row 11 $r.KellyFrac
    Set-TextCell $row 12 $r.Status
    Set-TextCell $row 13 ""
    Set-TextCell $row 14 ""
    Set-TextCell $row 15 ""
    Set-TextCell $row 16 ""
    Set-TextCell $row 17 ""
}
